# Update the LR-pairs (Ndp-Lgr4) sheet with the refreshed TPM-based results.
# The original sheet only had 4 data rows (MuSCs as the only sending cluster,
# against ECs/FAPs/MuSCs/Resolving-Mac target clusters). The refreshed script
# now reports both ECs and MuSCs as sending clusters, each against five target
# clusters (ECs, FAPs, Inflammatory-Mac, MuSCs, Resolving-Mac) -> 10 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Ndp", "Lgr4", "ECs", 1, 0.3333333333333333, 0.008202333333333334, 0.024607, 0.0618542966452333, 0.08999806156895879, 3, 1, 0.5779736666666667, 1.733921, 0.04835019606981441, 0.05356228267519154, 0.00474073267188889, 0.042666594047, 0.002990667370557494, 0.004820501613975863),
    @("ECs", "Ndp", "Lgr4", "FAPs", 1, 0.3333333333333333, 0.008202333333333334, 0.024607, 0.0618542966452333, 0.08999806156895879, 3, 1, 7.791016, 23.373048, 0.6517548686181108, 0.7220131735856595, 0.06390451023733335, 0.575140592136, 0.04031383898347968, 0.06497978604996152),
    @("ECs", "Ndp", "Lgr4", "Inflammatory-Mac", 1, 0.3333333333333333, 0.008202333333333334, 0.024607, 0.0618542966452333, 0.08999806156895879, 1, 0.3333333333333333, 0.05015166666666667, 0.150455, 0.004195421100317676, 0.004647681895481942, 0.0004113606872222223, 0.003702246185, 0.0002595048212907207, 0.0004182823613825189),
    @("ECs", "Ndp", "Lgr4", "MuSCs", 1, 0.3333333333333333, 0.008202333333333334, 0.024607, 0.0618542966452333, 0.08999806156895879, 2, 1, 3.489664, 6.979328, 0.2919266886169084, 0.215597330685123, 0.02862338734933333, 0.171740324096, 0.0180569199963709, 0.01940334184110287),
    @("ECs", "Ndp", "Lgr4", "Resolving-Mac", 1, 0.3333333333333333, 0.008202333333333334, 0.024607, 0.0618542966452333, 0.08999806156895879, 1, 0.3333333333333333, 0.0451, 0.1353, 0.003772825594848836, 0.004179531158543795, 0.0003699252333333334, 0.0033293271, 0.0002333654735345087, 0.0003761497025360061),
    @("MuSCs", "Ndp", "Lgr4", "ECs", 2, 1, 0.124405, 0.24881, 0.9381457033547667, 0.9100019384310413, 3, 1, 0.5779736666666667, 1.733921, 0.04835019606981441, 0.05356228267519154, 0.07190281400166666, 0.4314168840099999, 0.04535952869925691, 0.04874178106121568),
    @("MuSCs", "Ndp", "Lgr4", "FAPs", 2, 1, 0.124405, 0.24881, 0.9381457033547667, 0.9100019384310413, 3, 1, 7.791016, 23.373048, 0.6517548686181108, 0.7220131735856595, 0.9692413454799998, 5.81544807288, 0.6114410296346311, 0.6570333875356981),
    @("MuSCs", "Ndp", "Lgr4", "Inflammatory-Mac", 2, 1, 0.124405, 0.24881, 0.9381457033547667, 0.9100019384310413, 1, 0.3333333333333333, 0.05015166666666667, 0.150455, 0.004195421100317676, 0.004647681895481942, 0.006239118091666667, 0.03743470855, 0.003935916279026955, 0.004229399534099424),
    @("MuSCs", "Ndp", "Lgr4", "MuSCs", 2, 1, 0.124405, 0.24881, 0.9381457033547667, 0.9100019384310413, 2, 1, 3.489664, 6.979328, 0.2919266886169084, 0.215597330685123, 0.4341316499199999, 1.73652659968, 0.2738697686205375, 0.1961939888440202),
    @("MuSCs", "Ndp", "Lgr4", "Resolving-Mac", 2, 1, 0.124405, 0.24881, 0.9381457033547667, 0.9100019384310413, 1, 0.3333333333333333, 0.0451, 0.1353, 0.003772825594848836, 0.004179531158543795, 0.0056106655, 0.033663993, 0.003539460121314327, 0.003803381456007789)
)

$nrows = $data.Length
$ncols = $data[0].Length
$arr = New-Object 'object[,]' $nrows,$ncols
for ($i = 0; $i -lt $nrows; $i++) {
    for ($j = 0; $j -lt $ncols; $j++) {
        $arr[$i,$j] = $data[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $nrows - 1
$ws.Range("A$startRow`:T$endRow").Value = $arr
